$d = $word.ActiveDocument

# The document carries three inline "logo" pictures that were renamed
# (their Word display-name / <wp:docPr name="..."> attribute) without any
# other visible change: two copies of the Pearson Edexcel logo living in
# the page footers (image2.png -> image1.png) and the BTec logo living in
# the page header (image1.jpg -> image2.jpg). Walk every section's headers
# and footers and rename the pictures by matching on their (unique,
# unchanged) alternative text / description so the edit is robust
# regardless of section count.

foreach ($sec in $d.Sections) {
    foreach ($h in $sec.Headers) {
        if ($h.Exists -and $h.Range.InlineShapes.Count -gt 0) {
            foreach ($shp in $h.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
    foreach ($f in $sec.Footers) {
        if ($f.Exists -and $f.Range.InlineShapes.Count -gt 0) {
            foreach ($shp in $f.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
